# Adds a new "convert to number per kg" section (columns M/N) to sheet1,
# converting the hPa pressure inputs in columns I and L into number-per-kg
# units using the ideal gas law factor 101325/273.15/287, and sums the
# converted Case 1 / Case 2 values. Also bumps the sheet's zoom level.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New label heading the conversion columns
$ws.Range("M4").Value = "convert to number per kg"

# Case 1 (column I) and Case 2 (column L) pressure conversions
$ws.Range("M5").Formula = "=I5/(101325/273.15/287)"
$ws.Range("N5").Formula = "=L5/(101325/273.15/287)"

$ws.Range("M6").Formula = "=I6/(101325/273.15/287)"
$ws.Range("N6").Formula = "=L6/(101325/273.15/287)"

$ws.Range("M7").Formula = "=I7/(101325/273.15/287)"
$ws.Range("N7").Formula = "=L7/(101325/273.15/287)"

$ws.Range("N8").Formula = "=L8/(101325/273.15/287)"

# Totals
$ws.Range("M10").Formula = "=SUM(M5:M7)"
$ws.Range("N10").Formula = "=SUM(N5:N8)"

# Zoom the sheet view in
$excel.ActiveWindow.Zoom = 171
